$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the green "accent6" color from the "تاریخ ایجاد آگهی اضافه شود"
#    bullet (both the paragraph-mark run properties and the text run's
#    properties lose their <w:color> element; everything else is untouched).
# ---------------------------------------------------------------------------
$found = $d.Content.Find.Execute("تاریخ ایجاد آگهی اضافه شود", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the target paragraph text."
}

$target = $d.Content.Paragraphs.First
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*تاریخ ایجاد آگهی اضافه شود*") {
        $targetPara = $p
        break
    }
}
if ($null -eq $targetPara) {
    throw "Paragraph not found."
}

$fullRange = $targetPara.Range
$newParagraphXml = '<w:p w14:paraId="2921715B" w14:textId="328734C4" w:rsidR="00BB4287" w:rsidRPr="000C11C2" w:rsidRDefault="00F709AC" w:rsidP="000C11C2"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:bidi/><w:spacing w:line="480" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Vazir" w:hAnsi="Vazir" w:cs="Vazir"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr><w:r w:rsidRPr="000C11C2"><w:rPr><w:rFonts w:ascii="Vazir" w:hAnsi="Vazir" w:cs="Vazir" w:hint="cs"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>تاریخ ایجاد آگهی اضافه شود</w:t></w:r></w:p>'
$pkgXml = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml'><w:body>" + $newParagraphXml + "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
$fullRange.InsertXML($pkgXml)

# ---------------------------------------------------------------------------
# 2) Drop the two extra bullet points that were added under the "Edit" /
#    second-style list ("صفحه پروفایل کاربر از روی UI جدید زده شود" and
#    "تمام عکس هایی که دایره ای نشده ، دایره ای شود مثل صفحه پروفایل و من
#    همبرگری"). They are deleted completely, merging their neighbours back
#    together.
# ---------------------------------------------------------------------------
$toRemove = @(
    "صفحه پروفایل کاربر از روی UI جدید زده شود",
    "تمام عکس هایی که دایره ای نشده ، دایره ای شود مثل صفحه پروفایل و من همبرگری"
)

foreach ($needle in $toRemove) {
    $victim = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like ("*" + $needle + "*")) {
            $victim = $p
            break
        }
    }
    if ($null -eq $victim) {
        throw ("Paragraph to remove not found: " + $needle)
    }
    $victim.Range.Delete()
}

Write-Output "done"
